$d = $word.ActiveDocument

# --- 1) The three short-form occurrences: "...(Preview 14)" -> "...(RC1)" ---
#     - .NET MAUI App (Preview 14)
#     - .NET MAUI App (C#) (Preview 14)
#     - .NET MAUI Class Library (Preview 14)
# wdReplaceAll (11th arg = 2) so a single Execute call fixes all three.
$r = $d.Content
$r.Find.Execute("Preview 14)", $true, $false, $false, $false, $false, $true, 1, $false, "RC1)", 2)

# --- 2) The long-form occurrence in the closing "Note:" paragraph ---
#     ".NET MAUI Preview 14 (VS2022 17.2 Preview 2.0 or later)"
#  -> ".NET MAUI RC1 (VS2022 17.2 Preview 3.0 or later)"
$r = $d.Content
$r.Find.Execute("Preview 14 (VS2022 17.2 Preview 2.0 or later)", $true, $false, $false, $false, $false, $true, 1, $false, "RC1 (VS2022 17.2 Preview 3.0 or later)", 2)
